# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet (cloning the row1/row2 layout+style of the
#    "2021-Q4" sheet, which already uses the right header set/order) right
#    before the "总计" (totals) sheet, and fill in the new fund-holding row.
# 2. Prepend a "2022-Q1" row to the "总计" sheet's date/count/value table,
#    shifting the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q1" sheet
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheetBefore = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# Worksheet references returned by .Item(...) track a *position*, not a
# stable identity, so re-resolve "总计" by name now that the new sheet has
# been inserted in front of it (its index shifted from 4 to 5).
$totalSheet = $wb.Worksheets.Item("总计")

# Clone header row (B1:H1) + the A2 row-index cell from the template sheet so
# the new sheet picks up identical styling (bold/bordered header, etc.).
$template.Range("B1:H2").Copy($newSheet.Range("B1:H2"))
$template.Range("A2").Copy($newSheet.Range("A2"))

# Fund-code / ratios are stored as literal text (to keep the leading zero in
# the fund code and the exact decimal formatting), so force Text format
# before writing them.
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("D2:G2").NumberFormat = "@"

$newSheet.Range("B2").Value = "005585"
$newSheet.Range("C2").Value = "银河文体娱乐主题灵活配置混合"
$newSheet.Range("D2").Value = "5.54"
$newSheet.Range("E2").Value = "74.07"
$newSheet.Range("F2").Value = "4.11"
$newSheet.Range("G2").Value = "0.2277"
$newSheet.Range("H2").Value = 8

# ---------------------------------------------------------------------
# Step 2: prepend a "2022-Q1" row to the "总计" sheet
# ---------------------------------------------------------------------
$totalSheet.Rows(2).Insert()

# Pull the row-index cell style (bold/bordered, like A3:A5) onto the new A2,
# then clear the incidental formatting Insert() left on B2:D2 so they match
# the plain (unstyled) data cells used elsewhere in this table.
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.23

# Renumber the (0-based) row-index column for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
